$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Row 69: Emerald Mining Town
$ws.Cells.Item(69, 1).Value = "Emerald Mining Town"
$ws.Cells.Item(69, 2).Value = "Twisted Memories"
$ws.Cells.Item(69, 5).Value = "An old Dilapidated Mining town"
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 6
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 2352
$ws.Cells.Item(69, 11).Value = 64
$ws.Cells.Item(69, 13).Value = "Yes"

# Row 70: Twisted Memorial Crypt
$ws.Cells.Item(70, 1).Value = "Twisted Memorial Crypt"
$ws.Cells.Item(70, 2).Value = "Twisted Memories"
$ws.Cells.Item(70, 4).Value = "Twisted Memorial Crypt Key"
$ws.Cells.Item(70, 5).Value = "A crypt that contains the silence of darkness that was only once illuminated by the faint sound of a childs heart beat"
$ws.Cells.Item(70, 8).Value = 6
$ws.Cells.Item(70, 10).Value = 1760
$ws.Cells.Item(70, 11).Value = 384
$ws.Cells.Item(70, 13).Value = "Yes"
